$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Tue Oct 31 00:50:03 2023"
$ws.Range("B4").Value = "DJI"
$ws.Range("C4").Value = "Dow Jones Industrial Average"
$ws.Range("D4").Value = "Trading in Progress"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "32745.28"
$ws.Range("E4").Style = "Normal"

$ws.Range("F4").Value = "+327.69  +1.01%"
$ws.Range("G4").Value = 32802.36
$ws.Range("H4").Value = 32537.54
$ws.Range("I4").Value = 35679.13
$ws.Range("J4").Value = 0.73
$ws.Range("K4").Value = 32537.54
$ws.Range("L4").Value = 32417.59
$ws.Range("M4").Value = 31429.82
$ws.Range("N4").Value = 0.008
$ws.Range("O4").Value = 134000000
